# Fix methods to new data
# Updates computed summary statistics (mean/std regression results) on all
# four sheets to reflect re-run analysis results.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "species mean"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("species mean")
$ws.Range("C2").Value = 0.8989300411522642
$ws.Range("D2").Value = 6.252181069958859
$ws.Range("E2").Value = 0.5132522035361657
$ws.Range("F2").Value = 13.02461006533311
$ws.Range("G2").Value = 0.03831948025187215
$ws.Range("H2").Value = 499.3891749520004

$ws.Range("C3").Value = 2.566919191919193
$ws.Range("D3").Value = 7.480808080808081
$ws.Range("E3").Value = 1.045814780247639
$ws.Range("F3").Value = 48.00576749879822
$ws.Range("G3").Value = 0.03748718523446825
$ws.Range("H3").Value = 69.81733537295449

# ---------------------------------------------------------------------------
# Sheet "species bioshaker mean"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("species bioshaker mean")
$ws.Range("D2").Value = 1.565972222222224
$ws.Range("E2").Value = 6.772140522875823
$ws.Range("F2").Value = 0.6281495128170839
$ws.Range("G2").Value = 10.01068554655091
$ws.Range("H2").Value = 0.04290064197226771
$ws.Range("I2").Value = 667.0968282711668

$ws.Range("D3").Value = 0.2219320066334992
$ws.Range("E3").Value = 5.724461028192374
$ws.Range("F3").Value = 0.3966400090420997
$ws.Range("G3").Value = 16.08351853215685
$ws.Range("H3").Value = 0.03366994298341094
$ws.Range("I3").Value = 329.1784223295625

$ws.Range("D4").Value = 3.148842592592593
$ws.Range("E4").Value = 9.775000000000002
$ws.Range("F4").Value = 0.3597763375058016
$ws.Range("G4").Value = 22.87505852187092
$ws.Range("H4").Value = 0.04593606360508246
$ws.Range("I4").Value = 114.7833351028985

$ws.Range("D5").Value = 1.868611111111111
$ws.Range("E5").Value = 4.727777777777777
$ws.Range("F5").Value = 1.869060911537844
$ws.Range("G5").Value = 78.16261827111103
$ws.Range("H5").Value = 0.02734853118973121
$ws.Range("I5").Value = 15.85813569702169

# ---------------------------------------------------------------------------
# Sheet "Species std"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Species std")
$ws.Range("C2").Value = 0.8571820291471475
$ws.Range("D2").Value = 0.6504466319630029
$ws.Range("E2").Value = 0.139794714791673
$ws.Range("F2").Value = 3.778590780298534
$ws.Range("G2").Value = 0.005683596800217818
$ws.Range("H2").Value = 619.6949162757148

$ws.Range("C3").Value = 1.954582209721116
$ws.Range("D3").Value = 2.754640666345447
$ws.Range("E3").Value = 2.279423883375306
$ws.Range("F3").Value = 83.23337561277928
$ws.Range("G3").Value = 0.01354971319733675
$ws.Range("H3").Value = 81.68345397934863

# ---------------------------------------------------------------------------
# Sheet "Species bioshaker std"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Species bioshaker std")
$ws.Range("D2").Value = 0.4176650214687206
$ws.Range("E2").Value = 0.2130523963974698
$ws.Range("F2").Value = 0.04286243121554766
$ws.Range("G2").Value = 0.3043536084381681
$ws.Range("H2").Value = 0.001865346479553693
$ws.Range("I2").Value = 789.8194445996253

$ws.Range("D3").Value = 0.625300925371889
$ws.Range("E3").Value = 0.5016420471356491
$ws.Range("F3").Value = 0.102011886822773
$ws.Range("G3").Value = 3.168031339079627
$ws.Range("H3").Value = 0.00429935295426457
$ws.Range("I3").Value = 296.6905585105587

$ws.Range("D4").Value = 2.132991824232088
$ws.Range("F4").Value = 0.107722676439064
$ws.Range("G4").Value = 3.696566793870105
$ws.Range("H4").Value = 0.001633894002878522
$ws.Range("I4").Value = 87.95864472062848

$ws.Range("D5").Value = 1.536303958381193
$ws.Range("E5").Value = 1.505449064613414
$ws.Range("F5").Value = 3.27547223484301
$ws.Range("G5").Value = 119.564205707389
$ws.Range("H5").Value = 0.01468785579066686
$ws.Range("I5").Value = 13.46897201680261

$wb.Save()
